# Checkers Remote two player - Doug's updates for remote with repository connection
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move/record the active selection to M3 (was A1 by default)
$ws.Range("M3").Select()

# Fill in the "move counter" column (N) that was previously blank/stale,
# and add the matching "w/yellow" label in column O for row 5
$ws.Range("N5").Value = 5
$ws.Range("O5").Value = "w/yellow"

# Update the stale counter values in N6:N9 to the new sequential numbers
$ws.Range("N6").Value = 6
$ws.Range("N7").Value = 7
$ws.Range("N8").Value = 8
$ws.Range("N9").Value = 9
